$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)
$ftr = $sec.Footers.Item(1)

# Step 1: replace the old tail text with the new full tail text (this initially
# lands as a single merged run because the replacement text shares formatting
# with its neighbours).
$rng = $ftr.Range
$ok = $rng.Find.Execute(". 7 e 8 de dezembro de 2023. Câmara Municipal de Itapevi/SP", $true, $false, $false, $false, $false, $true, 1, $false, ". 18 e 19 de abril de 2024. Câmara Municipal de Itapevi/SP", 2)

# Absolute (story-relative) offset where the replacement text begins.
$base = $rng.Start

# Step 2: split the merged replacement text into the individual runs seen in the
# target document by re-asserting (toggling) formatting on each exact sub-range -
# this forces a run break at that boundary without altering the visible format.
$segments = @(". ", "18", " e ", "19", " de ", "abril ", "de 202", "4", ". Câmara Municipal de Itapevi/SP")

$pos = $base
foreach ($seg in $segments) {
    $piece = $ftr.Range
    $piece.Start = $pos
    $found = $piece.Find.Execute($seg, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    $piece.Bold = 1
    $piece.Bold = 0
    $pos = $pos + $seg.Length
}
